# Dev: start working on export module and finish data-export-json service.
#
# Adds a new "export" service-category block (5 rows) to the bottom of the
# services-details sheet, mirroring the existing "modify" block's formatting,
# and moves the sheet's viewport/selection down to the newly-added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clone formatting from the last existing group ("modify", row 45)
#        down onto the five new rows (46-50) before filling in values, so the
#        new block keeps the same alternating category fill/border style.
$ws.Range("A45:E45").Copy()
$ws.Range("A46:E50").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. Row numbers (column A)
$ws.Range("A46").Value = 48
$ws.Range("A47").Value = 49
$ws.Range("A48").Value = 50
$ws.Range("A49").Value = 51
$ws.Range("A50").Value = 52

# --- 3. Service names (column C) filled before the repeated "export"
#        category label (column B), matching the authoring order.
$ws.Range("C46").Value = "export-data-pdf"

$ws.Range("B46").Value = "export"
$ws.Range("B47").Value = "export"
$ws.Range("B48").Value = "export"
$ws.Range("B49").Value = "export"
$ws.Range("B50").Value = "export"

$ws.Range("C47").Value = "export-data-csv"
$ws.Range("C48").Value = "export-data-excel"
$ws.Range("C49").Value = "export-data-json"
$ws.Range("C50").Value = "export-data-xml"

# --- 4. Default ports (column D)
$ws.Range("D46").Value = 7030
$ws.Range("D47").Value = 7031
$ws.Range("D48").Value = 7032
$ws.Range("D49").Value = 7033
$ws.Range("D50").Value = 7034

# --- 5. Base paths (column E)
$ws.Range("E46").Value = "/selling/export/sales/pdf/"
$ws.Range("E47").Value = "/selling/export/sales/csv/"
$ws.Range("E48").Value = "/selling/export/sales/excel/"
$ws.Range("E49").Value = "/selling/export/sales/json/"
$ws.Range("E50").Value = "/selling/export/sales/xml/"

# --- 6. Leave the sheet scrolled/selected on the newly-added block, as the
#        author left it after typing the new rows.
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G39").Select()
